# live_trading_results.xlsx - append new/closed trades and refresh rollups
#
# Trade #101 (leadlag, UP) opened at 2026-02-16 21:40:15.
# Two previously-open momentum trades (#66 @ 21:35:06, #67 @ 21:35:12)
# closed out, so: momentum sheet rows flip OPEN -> CLOSED with fill
# price/P&L/exit-reason/duration, those two closed trades are appended
# to the "All Trades" ledger, and the Summary / Comparison roll-up
# numbers are refreshed to reflect the new totals.

$wb = $excel.ActiveWorkbook

# Helper: Excel's COM Value setter auto-detects numbers/dates/percentages
# from plain strings. Several cells in this workbook intentionally hold
# literal text that LOOKS like a number/date/percentage (e.g. "67.2%",
# "2026-02-16", "5.31"), so force those through as text with a leading
# apostrophe (exactly what typing '67.2% into Excel does) while leaving
# unambiguous text alone.
function Set-Text {
    param($range, [string]$text)
    if ($text -match '^[+-]?[0-9]') {
        $range.Value = "'" + $text
    } else {
        $range.Value = $text
    }
}

# ---------------------------------------------------------------------
# Summary sheet: OVERALL + momentum roll-up rows
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

$summary.Range("C2").Value = 67
Set-Text $summary.Range("D2") "67.2%"
Set-Text $summary.Range("E2") "+17.7744%"
Set-Text $summary.Range("F2") "+0.2653%"

$summary.Range("C4").Value = 25
Set-Text $summary.Range("D4") "48.0%"
Set-Text $summary.Range("E4") "+6.6831%"
Set-Text $summary.Range("F4") "+0.2673%"

# ---------------------------------------------------------------------
# leadlag sheet: append new open trade #101
# ---------------------------------------------------------------------
$leadlag = $wb.Worksheets.Item("leadlag")

$leadlag.Cells.Item(77, 1).Value = 101
Set-Text $leadlag.Cells.Item(77, 2) "2026-02-16"
Set-Text $leadlag.Cells.Item(77, 3) "21:40:15"
$leadlag.Cells.Item(77, 4).Value = "leadlag"
$leadlag.Cells.Item(77, 5).Value = "UP"
$leadlag.Cells.Item(77, 6).Value = 68362.545
$leadlag.Cells.Item(77, 8).Value = "OPEN"
$leadlag.Cells.Item(77, 9).Value = 0
$leadlag.Cells.Item(77, 10).Value = 0
$leadlag.Cells.Item(77, 11).Value = 0.75
$leadlag.Cells.Item(77, 12).Value = "Binance leading with 0.130% move"
$leadlag.Cells.Item(77, 14).Value = 0

# ---------------------------------------------------------------------
# momentum sheet: trades #66 (row 14) and #67 (row 15) close out
# ---------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")

$momentum.Cells.Item(14, 7).Value = 67765.526266
$momentum.Cells.Item(14, 8).Value = "CLOSED"
$momentum.Cells.Item(14, 9).Value = 1.0487
$momentum.Cells.Item(14, 10).Value = 10.49
$momentum.Cells.Item(14, 13).Value = "time_exit_5min"
$momentum.Cells.Item(14, 14).Value = 5

$momentum.Cells.Item(15, 7).Value = 67946.160194
$momentum.Cells.Item(15, 8).Value = "CLOSED"
$momentum.Cells.Item(15, 9).Value = 0.7836
$momentum.Cells.Item(15, 10).Value = 7.84
$momentum.Cells.Item(15, 13).Value = "time_exit_5min"
$momentum.Cells.Item(15, 14).Value = 5

# ---------------------------------------------------------------------
# All Trades sheet: append the two newly-closed momentum trades
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Cells.Item(67, 1).Value = 66
Set-Text $allTrades.Cells.Item(67, 2) "2026-02-16"
Set-Text $allTrades.Cells.Item(67, 3) "21:35:06"
$allTrades.Cells.Item(67, 4).Value = "momentum"
$allTrades.Cells.Item(67, 5).Value = "DOWN"
$allTrades.Cells.Item(67, 6).Value = 68483.69500000001
$allTrades.Cells.Item(67, 7).Value = 67765.526266
$allTrades.Cells.Item(67, 8).Value = "CLOSED"
$allTrades.Cells.Item(67, 9).Value = 1.0487
$allTrades.Cells.Item(67, 10).Value = 10.49
$allTrades.Cells.Item(67, 11).Value = 0.9
$allTrades.Cells.Item(67, 12).Value = "Downward momentum: -0.343% over 10 samples"
$allTrades.Cells.Item(67, 13).Value = "time_exit_5min"
$allTrades.Cells.Item(67, 14).Value = 5

$allTrades.Cells.Item(68, 1).Value = 67
Set-Text $allTrades.Cells.Item(68, 2) "2026-02-16"
Set-Text $allTrades.Cells.Item(68, 3) "21:35:12"
$allTrades.Cells.Item(68, 4).Value = "momentum"
$allTrades.Cells.Item(68, 5).Value = "DOWN"
$allTrades.Cells.Item(68, 6).Value = 68482.8
$allTrades.Cells.Item(68, 7).Value = 67946.160194
$allTrades.Cells.Item(68, 8).Value = "CLOSED"
$allTrades.Cells.Item(68, 9).Value = 0.7836
$allTrades.Cells.Item(68, 10).Value = 7.84
$allTrades.Cells.Item(68, 11).Value = 0.9
$allTrades.Cells.Item(68, 12).Value = "Downward momentum: -0.340% over 10 samples"
$allTrades.Cells.Item(68, 13).Value = "time_exit_5min"
$allTrades.Cells.Item(68, 14).Value = 5

# ---------------------------------------------------------------------
# Comparison sheet: momentum strategy roll-up row
# ---------------------------------------------------------------------
$comparison = $wb.Worksheets.Item("Comparison")

$comparison.Range("B3").Value = 25
Set-Text $comparison.Range("C3") "48.0%"
Set-Text $comparison.Range("D3") "6.94"
Set-Text $comparison.Range("E3") "+0.6506%"
Set-Text $comparison.Range("G3") "1.16"
